$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 665
$ws1.Range("G6").Value = "不可售"
$ws1.Range("F7").Value = 569
$ws1.Range("F10").Value = 1816
$ws1.Range("F11").Value = 1409
$ws1.Range("F13").Value = 1672
$ws1.Range("F17").Value = 22
$ws1.Range("F18").Value = 50
$ws1.Range("F26").Value = 94
$ws1.Range("F27").Value = 4206
$ws1.Range("F31").Value = 1998
$ws1.Range("F32").Value = 61
$ws1.Range("F33").Value = 1955

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G2").Value = 98

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 665
$ws4.Range("G6").Value = "不可售"
$ws4.Range("F7").Value = 569
$ws4.Range("G8").Value = 98
$ws4.Range("F11").Value = 1816
$ws4.Range("F12").Value = 1409
$ws4.Range("F14").Value = 1672
$ws4.Range("F18").Value = 22
$ws4.Range("F19").Value = 50
$ws4.Range("F27").Value = 94
$ws4.Range("F28").Value = 4206
$ws4.Range("F34").Value = 1998
$ws4.Range("F35").Value = 61
$ws4.Range("F36").Value = 1955
